# "Add files via upload" — re-upload of saitenuni.xlsx with:
#   1) two new empty worksheets (Sheet2, Sheet3) appended after Sheet1
#   2) a handful of score cells on Sheet1 bumped by +1 (and the two
#      running-total cells in row 93 updated to match)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Append Sheet2 then Sheet3 after the existing Sheet1, keeping tab order
# Sheet1, Sheet2, Sheet3 and restoring Sheet1 as the active/selected tab.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws1.Activate()

# Updated score values on Sheet1
$ws1.Range("C4").Value = 87
$ws1.Range("D4").Value = 68
$ws1.Range("C7").Value = 76
$ws1.Range("C24").Value = 70
$ws1.Range("D24").Value = 65
$ws1.Range("C36").Value = 41
$ws1.Range("D38").Value = 56
$ws1.Range("C46").Value = 71
$ws1.Range("D46").Value = 55
$ws1.Range("C63").Value = 107
$ws1.Range("D63").Value = 99
$ws1.Range("C68").Value = 61
$ws1.Range("D68").Value = 48
$ws1.Range("C92").Value = 250
$ws1.Range("D92").Value = 184

# Running totals (row 93) reflect the sum of the updated column values
$ws1.Range("C93").Value = 5513
$ws1.Range("D93").Value = 4445
